# Revert api endpoint validation / Fix race condition in concurrent processing.
# Appends a new row (row 28) of captured data to each of the four worksheets.

$wb = $excel.ActiveWorkbook

function Add-DataRow {
    param(
        [string]$SheetName,
        [string]$TimeStr,
        [string]$BVal,
        [string]$CVal,
        [string]$DVal,
        [string]$EVal,
        [double]$FVal,
        [string]$GStr,
        [double]$HVal,
        [double]$IVal
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $row = 28

    $timeVal = [double]$TimeStr
    $gVal = [double]$GStr

    $ws.Cells.Item($row, 1).Value = $timeVal
    $ws.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($row, 2).Value = $BVal
    $ws.Cells.Item($row, 3).Value = $CVal
    $ws.Cells.Item($row, 4).Value = $DVal
    $ws.Cells.Item($row, 5).Value = $EVal

    $ws.Cells.Item($row, 6).Value = $FVal
    $ws.Cells.Item($row, 7).Value = $gVal
    $ws.Cells.Item($row, 8).Value = $HVal
    $ws.Cells.Item($row, 9).Value = $IVal
}

Add-DataRow "ROW35-FE-LIFTER" "45738.29230510417" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c," "0x01,0x82" "0xd" 400 "568631262647113970876416" 386 13

Add-DataRow "ROW35-MID-LIFTER" "45738.14350274306" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c," "0x01,0x82" "0xe" 400 "568631262647113970876416" 386 14

Add-DataRow "ROW02-FE-LIFTER" "45738.29059998842" "0x01,0x90" "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c," "0x01,0x82" "0x3" 400 "568631262647113970876416" 386 3

Add-DataRow "ROW02-MID-LIFTER" "45738.34981295139" "0x01,0x90" "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c," "0x01,0x82" "0x3" 400 "985046333984776009023488" 386 3

$wb.Save()
